$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns D and E to Text format first so that numeric-looking
# price/percentage strings (e.g. "299.61") are kept as text, matching
# the inline-string cell content in the source data, instead of being
# auto-converted by Excel into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"
$ws.Range("E2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '42.777.97'
$ws.Range('D3').Value = '2.290.48'
$ws.Range('E3').Value = '  -0.82%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '299.61'
$ws.Range('E5').Value = '  -0.66%  '
$ws.Range('D6').Value = '96.49'
$ws.Range('E6').Value = '  -2.05%  '
$ws.Range('E7').Value = '  +0.83%  '
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').Value = '0.503'
$ws.Range('E9').Value = '  -3.58%  '
$ws.Range('D10').Value = '35.52'
$ws.Range('E10').Value = '  -0.38%  '
$ws.Range('E11').Value = '  -0.35%  '
$ws.Range('E12').Value = '  +0.70%  '
$ws.Range('E13').Value = '  -1.14%  '
$ws.Range('D14').Value = '6.73'
$ws.Range('D15').Value = '2.648.19'
$ws.Range('E15').Value = '  -0.81%  '
$ws.Range('D16').Value = '2.294.94'
$ws.Range('E16').Value = '  -0.35%  '
$ws.Range('D17').Value = '0.773'
$ws.Range('E17').Value = '  -1.93%  '
$ws.Range('D18').Value = '42.712.51'
$ws.Range('E18').Value = '  -0.62%  '
$ws.Range('E19').Value = '  -4.64%  '
$ws.Range('E20').Value = '  -0.47%  '
$ws.Range('D21').Value = '6.04'
$ws.Range('E21').Value = '  -2.21%  '
$ws.Range('D22').Value = '67.69'
$ws.Range('D23').Value = '240.51'
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('E24').Value = '  -1.99%  '
$ws.Range('E25').Value = '  +0.09%  '
$ws.Range('B26').Value = 'LEO'
$ws.Range('C26').Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D26').Value = '4.02'
$ws.Range('E26').Value = '  -0.22%  '
$ws.Range('B27').Value = 'PancakeSwap'
$ws.Range('C27').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D27').Value = '2.41'
$ws.Range('E27').Value = '  -1.19%  '
$ws.Range('D28').Value = '25.08'
$ws.Range('E28').Value = '  +0.44%  '
$ws.Range('D29').Value = '165.69'
$ws.Range('E29').Value = '  -1.96%  '
$ws.Range('E30').Value = '  -1.45%  '
$ws.Range('D31').Value = '9.01'
$ws.Range('E31').Value = '  -1.56%  '
$ws.Range('D32').Value = '32.72'
$ws.Range('E32').Value = '  -1.90%  '
$ws.Range('E33').Value = '  +0.06%  '
$ws.Range('D34').Value = '4.82'
$ws.Range('E34').Value = '  -1.20%  '
$ws.Range('D35').Value = '4.99'
$ws.Range('E35').Value = '  -4.07%  '
$ws.Range('D36').Value = '16.86'
$ws.Range('E36').Value = '  -7.83%  '
$ws.Range('E37').Value = '  -1.74%  '
$ws.Range('E38').Value = '  -1.40%  '
$ws.Range('E39').Value = '  -1.50%  '
$ws.Range('D40').Value = '1.74'
$ws.Range('E40').Value = '  -3.49%  '
$ws.Range('E41').Value = '  +0.16%  '
$ws.Range('D42').Value = '2.72'
$ws.Range('E42').Value = '  -0.98%  '
$ws.Range('D43').Value = '2.011.65'
$ws.Range('E43').Value = '  +1.10%  '
$ws.Range('E44').Value = '  -2.66%  '
$ws.Range('D45').Value = '10.05'
$ws.Range('E45').Value = '  -0.58%  '
$ws.Range('D46').Value = '2.10'
$ws.Range('E46').Value = '  +1.62%  '
$ws.Range('D47').Value = '17.12'
$ws.Range('E47').Value = '  -1.97%  '
$ws.Range('D48').Value = '2.77'
$ws.Range('E48').Value = '  -2.16%  '
$ws.Range('D49').Value = '2.87'
$ws.Range('E49').Value = '  -3.70%  '
$ws.Range('D50').Value = '2.515.50'
$ws.Range('E50').Value = '  -0.91%  '
$ws.Range('D51').Value = '52.87'
$ws.Range('E51').Value = '  -3.10%  '

# Restore the default (Normal) style on these columns so the cells
# do not retain an explicit text-number-format style attribute.
$ws.Range("D2:D51").Style = "Normal"
$ws.Range("E2:E51").Style = "Normal"
